$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = "[Shao-Fang%Nie%NULL%2,                   Miao%Yu%NULL%2,                   Tian%Xie%NULL%2,                   Fen%Yang%NULL%2,                   Hong-Bo%Wang%NULL%2,                   Zhao-Hui%Wang%NULL%2,                   Ming%Li%NULL%4,                   Xing-Li%Gao%NULL%2,                   Bing-Jie%Lv%NULL%2,                   Shi-Jia%Wang%NULL%2,                   Xiao-Bo%Zhang%NULL%2,                   Shao-Lin%He%NULL%2,                   Zhi-Hua%Qiu%NULL%2,                   Yu-Hua%Liao%NULL%2,                   Zi-Hua%Zhou%NULL%2,                   Xiang%Cheng%NULL%2]"
$ws.Cells.Item(2, 9).Value = "_PMC_CROSSREF"
$ws.Cells.Item(3, 5).Value = "[Christopher M%Petrilli%NULL%0,                   Simon A%Jones%NULL%3,                   Jie%Yang%NULL%5,                   Harish%Rajagopalan%NULL%3,                   Luke%O’Donnell%NULL%3,                   Yelena%Chernyak%NULL%3,                   Katie A%Tobin%NULL%3,                   Robert J%Cerfolio%NULL%3,                   Fritz%Francois%NULL%3,                   Leora I%Horwitz%NULL%3]"
$ws.Cells.Item(3, 9).Value = "_PMC_CROSSREF"
$ws.Cells.Item(4, 5).Value = "[Eboni G.%Price-Haywood%NULL%0,                   Jeffrey%Burton%NULL%2,                   Jeffrey%Burton%NULL%0,                   Daniel%Fort%NULL%1,                   Leonardo%Seoane%NULL%1]"
$ws.Cells.Item(4, 9).Value = "_PMC_CROSSREF"
$ws.Cells.Item(5, 5).Value = "[Mohamad%Raad%NULL%1,                   Mohammed%Dabbagh%NULL%1,                   Sarah%Gorgis%NULL%1,                   Jerry%Yan%NULL%1,                   Omar%Chehab%NULL%1,                   Carina%Dagher%NULL%2,                   Khaled%Jamoor%NULL%1,                   Inaya Hajj%Hussein%NULL%1,                   Bernard%Cook%NULL%1,                   Meredith%Van Harn%NULL%1,                   Gurjit%Singh%NULL%1,                   James%McCord%NULL%1,                   Sachin%Parikh%NULL%1]"
$ws.Cells.Item(5, 9).Value = "_PMC_elsevier_CROSSREF"
$ws.Cells.Item(6, 5).Value = "[Priyank%Shah%NULL%0,                   Rajkumar%Doshi%NULL%1,                   Avantika%Chenna%NULL%1,                   Robin%Owens%NULL%1,                   Abigail%Cobb%NULL%1,                   Holley%Ivey%NULL%1,                   Sarah%Newton%NULL%1,                   Kelly%Mccarley%NULL%1]"
$ws.Cells.Item(6, 9).Value = "_PMC_elsevier_CROSSREF"
$ws.Cells.Item(7, 5).Value = "[Sachin J.%Shah%sachin.shah@ucsf.edu%1,                   Peter N.%Barish%NULL%1,                   Priya A.%Prasad%NULL%1,                   Amy%Kistler%NULL%1,                   Norma%Neff%NULL%1,                   Jack%Kamm%NULL%1,                   Lucy M.%Li%NULL%1,                   Charles Y.%Chiu%NULL%1,                   Jennifer M.%Babik%NULL%1,                   Margaret C.%Fang%NULL%1,                   Yumiko%Abe-Jones%NULL%1,                   Narges%Alipanah%NULL%1,                   Francisco N.%Alvarez%NULL%1,                   Olga Borisovna%Botvinnik%NULL%1,                   Gloria%Castaneda%NULL%1,                   NULL%NULL%NULL%0,                   Rand M.%Dadasovich%NULL%1,                   Jennifer%Davis%NULL%1,                   Xianding%Deng%NULL%1,                   Joseph L.%DeRisi%NULL%1,                   Angela M.%Detweiler%NULL%1,                   Scot%Federman%NULL%1,                   John%Haliburton%NULL%1,                   Samantha%Hao%NULL%1,                   Andrew D.%Kerkhoff%NULL%1,                   G. Renuka%Kumar%NULL%1,                   Katherine B.%Malcolm%NULL%1,                   Sabrina A.%Mann%NULL%1,                   Sandra%Martinez%NULL%1,                   Rupa K.%Mary%NULL%1,                   Eran%Mick%NULL%1,                   Lusajo%Mwakibete%NULL%1,                   Nader%Najafi%NULL%1,                   Michael J.%Peluso%NULL%1,                   Maira%Phelps%NULL%1,                   Angela Oliveira%Pisco%NULL%1,                   Kalani%Ratnasiri%NULL%1,                   Luis A.%Rubio%NULL%1,                   Anna%Sellas%NULL%1,                   Kyla D.%Sherwood%NULL%1,                   Jonathan%Sheu%NULL%1,                   Natasha%Spottiswoode%NULL%1,                   Michelle%Tan%NULL%1,                   Guixia%Yu%NULL%1,                   Kirsten Neudoerffer%Kangelaris%NULL%1,                   Charles%Langelier%chaz.langelier@ucsf.edu%1]"
$ws.Cells.Item(7, 9).Value = "_PMC_elsevier_CROSSREF"
$ws.Cells.Item(8, 5).Value = "[Yinzhong%Shen%NULL%0,                   Fang%Zheng%NULL%3,                   Fang%Zheng%NULL%0,                   Danfeng%Sun%NULL%1,                   Yun%Ling%NULL%3,                   Jun%Chen%NULL%4,                   Feng%Li%NULL%3,                   Feng%Li%NULL%0,                   Tao%Li%NULL%0,                   Zhiping%Qian%NULL%1,                   Yuyi%Zhang%NULL%1,                   Qingnian%Xu%NULL%1,                   Li%Liu%NULL%0,                   Qin%Huang%NULL%1,                   Fei%Shan%NULL%1,                   Lie%Xu%NULL%1,                   Jun%Wu%NULL%1,                   Zhaoqin%Zhu%NULL%1,                   Zhigang%Song%NULL%1,                   Shenyang%Li%NULL%1,                   Yuxin%Shi%NULL%1,                   Jianliang%Zhang%NULL%1,                   Xueyun%Wu%NULL%1,                   Joshua B.%Mendelsohn%NULL%2,                   Joshua B.%Mendelsohn%NULL%0,                   Tongyu%Zhu%NULL%1,                   Hongzhou%Lu%NULL%3]"
$ws.Cells.Item(8, 9).Value = "_PMC_CROSSREF"
$ws.Cells.Item(9, 9).Value = "_CROSSREF"
$ws.Cells.Item(10, 5).Value = "[F. M. A.%van den Heuvel%NULL%1,                   J. L.%Vos%NULL%1,                   Y.%Koop%NULL%1,                   A. P. J.%van Dijk%NULL%1,                   A. L.%Duijnhouwer%NULL%1,                   Q.%de Mast%NULL%1,                   F. L.%van de Veerdonk%NULL%1,                   F.%Bosch%NULL%1,                   B.%Kok%NULL%1,                   M. G.%Netea%NULL%1,                   J.%Hoogerwerf%NULL%1,                   W.%Hoefsloot%NULL%1,                   E. T. T. L.%Tjwa%NULL%1,                   C. L.%de Korte%NULL%1,                   R. R. J.%van Kimmenade%NULL%1,                   R.%Nijveldt%robin@nijveldt.net%1]"
$ws.Cells.Item(10, 9).Value = "_PMC_Springer_CROSSREF"
$ws.Cells.Item(11, 5).Value = "[Jia-Fu%Wei%NULL%2,                   Fang-Yang%Huang%NULL%2,                   Tian-Yuan%Xiong%NULL%4,                   Tian-Yuan%Xiong%NULL%0,                   Qi%Liu%NULL%2,                   Hong%Chen%NULL%3,                   Hui%Wang%NULL%0,                   He%Huang%NULL%3,                   Yi-Chun%Luo%NULL%2,                   Xuan%Zhou%NULL%2,                   Zhi-Yue%Liu%NULL%2,                   Yong%Peng%NULL%2,                   Yuan-Ning%Xu%NULL%2,                   Bo%Wang%NULL%0,                   Ying-Ying%Yang%NULL%2,                   Zong-An%Liang%NULL%3,                   Xue-Zhong%Lei%NULL%2,                   Yang%Ge%NULL%2,                   Ming%Yang%NULL%3,                   Ling%Zhang%NULL%0,                   Ming-Quan%Zeng%NULL%2,                   He%Yu%NULL%2,                   Kai%Liu%NULL%2,                   Yu-Heng%Jia%NULL%2,                   Bernard D%Prendergast%NULL%2,                   Wei-Min%Li%NULL%5,                   Wei-Min%Li%NULL%0,                   Mao%Chen%NULL%2]"
$ws.Cells.Item(11, 9).Value = "_PMC_CROSSREF"
$ws.Cells.Item(12, 5).Value = "[Yuan%Yu%NULL%4,                   Dan%Xu%NULL%5,                   Shouzhi%Fu%NULL%5,                   Jun%Zhang%NULL%4,                   Xiaobo%Yang%NULL%4,                   Liang%Xu%NULL%4,                   Jiqian%Xu%NULL%9,                   Yongran%Wu%NULL%9,                   Chaolin%Huang%NULL%4,                   Yaqi%Ouyang%NULL%4,                   Luyu%Yang%NULL%4,                   Minghao%Fang%NULL%11,                   Hongwen%Xiao%NULL%4,                   Jing%Ma%NULL%4,                   Wei%Zhu%NULL%0,                   Song%Hu%NULL%4,                   Quan%Hu%NULL%4,                   Daoyin%Ding%NULL%4,                   Ming%Hu%NULL%8,                   Guochao%Zhu%NULL%4,                   Weijiang%Xu%NULL%4,                   Jun%Guo%NULL%4,                   Jinglong%Xu%NULL%4,                   Haitao%Yuan%NULL%4,                   Bin%Zhang%NULL%5,                   Zhui%Yu%yuzhui@whu.edu.cn%9,                   Dechang%Chen%icudechangchen@163.com%4,                   Shiying%Yuan%yuan_shiying@163.com%9,                   You%Shang%you_shanghust@163.com%10]"
$ws.Cells.Item(12, 9).Value = "_PMC_Springer_CROSSREF"
$ws.Cells.Item(13, 5).Value = "[Jia-Hui%Zeng%zengjiahui1993@163.com%1,                   Wei-Bo%Wu%13923451469@139.com%2,                   Wei-Bo%Wu%13923451469@139.com%0,                   Jiu-Xin%Qu%qujiuxin@163.com%1,                   Yao%Wang%68848030@qq.com%1,                   Chang-Feng%Dong%dongchangfeng7@163.com%1,                   Yong-Fang%Luo%luoyongfang2005@foxmail.com%1,                   Dan%Zhou%506057507@qq.com%1,                   Wen-Xia%Feng%1484177993@qq.com%1,                   Cheng%Feng%chaosheng-01@szsy.sustech.edu.cn%1]"
$ws.Cells.Item(13, 9).Value = "_PMC_Springer_CROSSREF"
$ws.Cells.Item(14, 5).Value = "[Xin-Ying%Zhao%NULL%1,                   Xuan-Xuan%Xu%NULL%1,                   Hai-Sen%Yin%NULL%1,                   Qin-Ming%Hu%NULL%1,                   Tao%Xiong%NULL%1,                   Yuan-Yan%Tang%NULL%1,                   Ai-Ying%Yang%NULL%1,                   Bao-Ping%Yu%yubp62@163.com%1,                   Zhi-Ping%Huang%huangzhipingdoctor@163.com%1]"
$ws.Cells.Item(14, 9).Value = "_PMC_Springer_CROSSREF"
$ws.Cells.Item(15, 5).Value = "[Fei%Zhou%NULL%0,                   Ting%Yu%NULL%0,                   Ronghui%Du%NULL%0,                   Guohui%Fan%NULL%0,                   Ying%Liu%NULL%0,                   Zhibo%Liu%NULL%0,                   Jie%Xiang%NULL%0,                   Yeming%Wang%NULL%0,                   Bin%Song%NULL%0,                   Xiaoying%Gu%NULL%0,                   Lulu%Guan%NULL%0,                   Yuan%Wei%NULL%0,                   Hui%Li%NULL%0,                   Xudong%Wu%NULL%0,                   Jiuyang%Xu%NULL%0,                   Shengjin%Tu%NULL%0,                   Yi%Zhang%NULL%0,                   Hua%Chen%NULL%0,                   Bin%Cao%NULL%0]"
$ws.Cells.Item(15, 9).Value = "_PMC_elsevier_CROSSREF"
$ws.Cells.Item(16, 5).Value = "[Aquino Bruno%Heberto%NULL%1,                   Plata Corona Juan%Carlos%vic_plata@hotmail.com%1,                   Castro Rubio José%Antonio%NULL%1,                   Pulido Pérez%Patricia%NULL%1,                   Torres Rasgado%Enrique%NULL%1,                   Morales Portano Julieta%Danira%NULL%1,                   Gómez Álvarez Enrique%Benito%NULL%1,                   Merino Rajme José%Alfredo%NULL%1]"
$ws.Cells.Item(16, 9).Value = "_PMC_elsevier_CROSSREF"
$ws.Cells.Item(17, 5).Value = "[Ling%Hu%NULL%2,                   Shaoqiu%Chen%NULL%1,                   Yuanyuan%Fu%NULL%1,                   Zitong%Gao%NULL%1,                   Hui%Long%NULL%1,                   Jian-ming%Wang%NULL%0,                   Hong-wei%Ren%NULL%1,                   Yi%Zuo%NULL%1,                   Huan%Li%NULL%1,                   Jie%Wang%NULL%1,                   Qing-bang%Xu%NULL%1,                   Wen-xiong%Yu%NULL%1,                   Jia%Liu%NULL%1,                   Chen%Shao%NULL%1,                   Jun-jie%Hao%NULL%1,                   Chuan-zhen%Wang%NULL%1,                   Yao%Ma%NULL%0,                   Zhanwei%Wang%NULL%1,                   Richard%Yanagihara%NULL%1,                   Youping%Deng%dengy@hawaii.edu%1]"
$ws.Cells.Item(17, 9).Value = "_PMC_CROSSREF"
$ws.Cells.Item(18, 5).Value = "[Mao%Huang%NULL%0,                   Yi%Yang%NULL%1,                   Futai%Shang%NULL%1,                   Yishan%Zheng%NULL%1,                   Wenjing%Zhao%NULL%1,                   Liang%Luo%NULL%1,                   Xudong%Han%NULL%1,                   Aihua%Lin%NULL%1,                   Hongsheng%Zhao%NULL%1,                   Qing%Gu%NULL%1,                   Yi%Shi%NULL%1,                   Jun%Li%NULL%1,                   Xingxiang%Xu%NULL%1,                   Kexi%Liu%NULL%1,                   YiJun%Deng%NULL%1,                   Quan%Cao%NULL%1,                   Weiwei%Wang%NULL%1]"
$ws.Cells.Item(18, 9).Value = "_PMC_elsevier_CROSSREF"
$ws.Cells.Item(19, 5).Value = "[Chiara%Lazzeri%NULL%1,                   Manuela%Bonizzoli%NULL%1,                   Stefano%Batacchi%NULL%1,                   Giovanni%Cianchi%NULL%1,                   Andrea%Franci%NULL%1,                   Giorgio Enzo%Fulceri%NULL%1,                   Adriano%Peris%NULL%1]"
$ws.Cells.Item(19, 9).Value = "_PMC_elsevier_CROSSREF"
$ws.Cells.Item(20, 5).Value = "[Chenze%Li%NULL%2,                   Jiangang%Jiang%NULL%3,                   Feng%Wang%NULL%0,                   Ning%Zhou%NULL%6,                   Giacomo%Veronese%NULL%2,                   Javid J.%Moslehi%NULL%3,                   Enrico%Ammirati%NULL%3,                   Dao Wen%Wang%NULL%3]"
$ws.Cells.Item(20, 9).Value = "_PMC_elsevier_CROSSREF"
$ws.Cells.Item(21, 5).Value = "[Alvaro%Lorente-Ros%xref no email%1,            Juan Manuel%Monteagudo Ruiz%xref no email%1,            Luis M.%Rinc\u00f3n%xref no email%1,            Rodrigo%Ortega P\u00e9rez%xref no email%1,            Sonia%Rivas%xref no email%1,            Rafael%Mart\u00ednez-Moya%xref no email%1,            Maria Ascensi\u00f3n%Sanrom\u00e1n%xref no email%1,            Luis%Manzano%xref no email%1,            Gonzalo Luis%Alonso%xref no email%1,            Borja%Ib\u00e1\u00f1ez%xref no email%1,            Jose Luis%Zamorano%xref no email%1]"
$ws.Cells.Item(21, 9).Value = "_CROSSREF"
$ws.Cells.Item(22, 5).Value = "[David T.%Majure%NULL%1,                   Luis%Gruberg%NULL%1,                   Shahryar G.%Saba%NULL%1,                   Charlotte%Kvasnovsky%NULL%1,                   Jamie S.%Hirsch%NULL%1,                   Rajiv%Jauhar%NULL%1,                   NULL%NULL%NULL%0]"
$ws.Cells.Item(22, 9).Value = "_PMC_elsevier_CROSSREF"
$ws.Cells.Item(23, 5).Value = "[Rong-Hui%Du%NULL%0,                   Li-Rong%Liang%NULL%2,                   Cheng-Qing%Yang%NULL%2,                   Wen%Wang%NULL%2,                   Tan-Ze%Cao%NULL%2,                   Ming%Li%NULL%0,                   Guang-Yun%Guo%NULL%2,                   Juan%Du%NULL%2,                   Chun-Lan%Zheng%NULL%2,                   Qi%Zhu%NULL%2,                   Ming%Hu%NULL%0,                   Xu-Yan%Li%NULL%2,                   Peng%Peng%NULL%0,                   Huan-Zhong%Shi%NULL%2]"
$ws.Cells.Item(23, 9).Value = "_PMC_CROSSREF"
$ws.Cells.Item(24, 5).Value = "[Jessica%Ferguson%NULL%2,                   Joelle I.%Rosser%NULL%2,                   Orlando%Quintero%NULL%2,                   Jake%Scott%NULL%2,                   Aruna%Subramanian%NULL%2,                   Mohammad%Gumma%NULL%2,                   Angela%Rogers%NULL%2,                   Shanthi%Kappagoda%NULL%2]"
$ws.Cells.Item(24, 9).Value = "_PMC_CROSSREF"
$ws.Cells.Item(25, 5).Value = "[Caroline E%Franks%NULL%1,                   Mitchell G%Scott%mgscott@wustl.edu%1,                   Christopher W%Farnsworth%NULL%1]"
$ws.Cells.Item(25, 9).Value = "_PMC_CROSSREF"
$ws.Cells.Item(26, 5).Value = "[Michael%Gottlieb%michaelgottliebmd@gmail.com%1,                   Sarah%Sansom%NULL%2,                   Sarah%Sansom%NULL%0,                   Casey%Frankenberger%NULL%1,                   Edward%Ward%NULL%1,                   Bala%Hota%NULL%1,                   Timothy%Jang%NULL%2,                   Timothy%Jang%NULL%0]"
$ws.Cells.Item(26, 9).Value = "_PMC_CROSSREF"
$ws.Cells.Item(27, 5).Value = "[Parag%Goyal%NULL%4,                   Justin J.%Choi%NULL%3,                   Laura C.%Pinheiro%NULL%3,                   Edward J.%Schenck%NULL%3,                   Ruijun%Chen%NULL%5,                   Assem%Jabri%NULL%3,                   Michael J.%Satlin%NULL%3,                   Thomas R.%Campion%NULL%3,                   Musarrat%Nahid%NULL%3,                   Joanna B.%Ringel%NULL%3,                   Katherine L.%Hoffman%NULL%3,                   Mark N.%Alshak%NULL%3,                   Han A.%Li%NULL%3,                   Graham T.%Wehmeyer%NULL%6,                   Graham T.%Wehmeyer%NULL%0,                   Mangala%Rajan%NULL%3,                   Evgeniya%Reshetnyak%NULL%3,                   Nathaniel%Hupert%NULL%3,                   Evelyn M.%Horn%NULL%3,                   Fernando J.%Martinez%NULL%3,                   Roy M.%Gulick%NULL%3,                   Monika M.%Safford%NULL%6,                   Monika M.%Safford%NULL%0]"
$ws.Cells.Item(27, 9).Value = "_PMC_CROSSREF"
$ws.Cells.Item(28, 5).Value = "[Anjali%Bhatla%NULL%1,                   Michael M.%Mayer%NULL%1,                   Srinath%Adusumalli%NULL%1,                   Matthew C.%Hyman%NULL%1,                   Eric%Oh%NULL%1,                   Ann%Tierney%NULL%1,                   Juwann%Moss%NULL%1,                   Anwar A.%Chahal%NULL%1,                   George%Anesi%NULL%1,                   Srinivas%Denduluri%NULL%1,                   Christopher M.%Domenico%NULL%1,                   Jeffrey%Arkles%NULL%1,                   Benjamin S.%Abella%NULL%1,                   John R.%Bullinga%NULL%1,                   David J.%Callans%NULL%1,                   Sanjay%Dixit%NULL%1,                   Andrew E.%Epstein%NULL%1,                   David S.%Frankel%NULL%1,                   Fermin C.%Garcia%NULL%1,                   Ramanan%Kumareswaram%NULL%1,                   Saman%Nazarian%NULL%1,                   Michael P.%Riley%NULL%1,                   Pasquale%Santangeli%NULL%1,                   Robert D.%Schaller%NULL%1,                   Gregory E.%Supple%NULL%1,                   David%Lin%NULL%1,                   Francis%Marchlinski%NULL%1,                   Rajat%Deo%Rajat.Deo@pennmedicine.upenn.edu%1]"
$ws.Cells.Item(28, 9).Value = "_PMC_elsevier_CROSSREF"
$ws.Cells.Item(29, 5).Value = "[Pavan K.%Bhatraju%NULL%0,                   Bijan J.%Ghassemieh%NULL%0,                   Michelle%Nichols%NULL%0,                   Richard%Kim%NULL%0,                   Keith R.%Jerome%NULL%0,                   Arun K.%Nalla%NULL%0,                   Alexander L.%Greninger%NULL%0,                   Sudhakar%Pipavath%NULL%0,                   Mark M.%Wurfel%NULL%0,                   Laura%Evans%NULL%0,                   Patricia A.%Kritek%NULL%0,                   T. Eoin%West%NULL%0,                   Andrew%Luks%NULL%0,                   Anthony%Gerbino%NULL%0,                   Chris R.%Dale%NULL%0,                   Jason D.%Goldman%NULL%0,                   Shane%O’Mahony%NULL%0,                   Carmen%Mikacenic%NULL%0]"
$ws.Cells.Item(29, 9).Value = "_PMC_CROSSREF"
$ws.Cells.Item(30, 5).Value = "[Frederick S%Buckner%fbuckner@uw.edu%4,                   Denise J%McCulloch%NULL%4,                   Denise J%McCulloch%NULL%0,                   Vidya%Atluri%NULL%2,                   Michela%Blain%NULL%2,                   Sarah A%McGuffin%NULL%2,                   Arun K%Nalla%NULL%2,                   Meei-Li%Huang%NULL%2,                   Alex L%Greninger%NULL%2,                   Keith R%Jerome%NULL%2,                   Seth A%Cohen%NULL%2,                   Santiago%Neme%NULL%2,                   Margaret L%Green%NULL%2,                   Helen Y%Chu%NULL%2,                   H Nina%Kim%NULL%2]"
$ws.Cells.Item(30, 9).Value = "_PMC_CROSSREF"
$ws.Cells.Item(31, 5).Value = "[Luca%Arcari%luca_arcari@outlook.it%1,                   Michelangelo%Luciani%NULL%2,                   Michelangelo%Luciani%NULL%0,                   Luca%Cacciotti%NULL%1,                   Maria Beatrice%Musumeci%NULL%1,                   Valerio%Spuntarelli%NULL%1,                   Eleonora%Pistella%NULL%1,                   Dario%Martolini%NULL%1,                   Daniele%Manzo%NULL%1,                   Mariateresa%Pucci%NULL%1,                   Claudio%Marone%NULL%1,                   Serena%Melandri%NULL%1,                   Gerardo%Ansalone%NULL%1,                   Claudio%Santini%NULL%1,                   Paolo%Martelletti%NULL%1,                   Massimo%Volpe%NULL%1,                   Luciano%De Biase%NULL%1]"
$ws.Cells.Item(31, 9).Value = "_PMC_Springer_CROSSREF"
$ws.Cells.Item(32, 5).Value = "[Elie%Azoulay%elie.azoulay@aphp.fr%1,                   Muriel%Fartoukh%NULL%2,                   Muriel%Fartoukh%NULL%0,                   Michael%Darmon%NULL%1,                   Guillaume%Géri%NULL%1,                   Guillaume%Voiriot%NULL%1,                   Thibault%Dupont%NULL%1,                   Lara%Zafrani%NULL%1,                   Lola%Girodias%NULL%1,                   Vincent%Labbé%NULL%1,                   Martin%Dres%NULL%1,                   Alexandra%Beurton%NULL%1,                   Antoine%Vieillard-Baron%NULL%1,                   Alexandre%Demoule%NULL%1]"
$ws.Cells.Item(32, 9).Value = "_PMC_Springer_CROSSREF"
$ws.Cells.Item(33, 5).Value = "[Hasan Ali%Barman%NULL%1,                   Adem%Atici%NULL%1,                   Irfan%Sahin%NULL%1,                   Gokhan%Alici%NULL%1,                   Esra%Aktas Tekin%NULL%1,                   Ömer Faruk%Baycan%NULL%1,                   Fatih%Ozturk%NULL%1,                   Ersan%Oflar%NULL%1,                   Sevil%Tugrul%NULL%1,                   Mustafa baran%Yavuz%NULL%1,                   Fatma Betul%Celik%NULL%1,                   Aysu%Oktay%NULL%1,                   Haluk%Vahaboglu%NULL%1,                   Mine%Adas%NULL%1,                   Namigar%Turgut%NULL%1,                   Ertugrul%Okuyan%NULL%1,                   Mustafa Taner%Yildirmak%NULL%1,                   Baris%Gungor%NULL%1]"
$ws.Cells.Item(33, 9).Value = "_PMC_CROSSREF"
$ws.Cells.Item(34, 5).Value = "[Alberto%Cipriani%alberto.cipriani@unipd.it%1,                   Federico%Capone%NULL%2,                   Federico%Capone%NULL%0,                   Filippo%Donato%NULL%1,                   Leonardo%Molinari%NULL%1,                   Davide%Ceccato%NULL%1,                   Alois%Saller%NULL%1,                   Lorenzo%Previato%NULL%1,                   Raffaele%Pesavento%NULL%1,                   Cristiano%Sarais%NULL%1,                   Paola%Fioretto%NULL%2,                   Sabino%Iliceto%NULL%1,                   Dario%Gregori%NULL%1,                   Angelo%Avogaro%NULL%1,                   Roberto%Vettor%NULL%3]"
$ws.Cells.Item(34, 9).Value = "_PMC_Springer_CROSSREF"
$ws.Cells.Item(35, 5).Value = "[Cao%Y%coreGivesNoEmail%0,                  Imam%Z%coreGivesNoEmail%1,                  Lippi%G%coreGivesNoEmail%1,                  Oran%DP%coreGivesNoEmail%1,                  Shi%S%coreGivesNoEmail%1]"
$ws.Cells.Item(35, 9).Value = "_CROSSREF"
$ws.Cells.Item(36, 5).Value = "[Juan M%Mej\u00eda-Vilet%xref no email%1,            Bertha M.%C\u00f3rdova-S\u00e1nchez%xref no email%1,            Dheni A.%Fern\u00e1ndez-Camargo%xref no email%1,            R. Ang\u00e9lica%M\u00e9ndez-P\u00e9rez%xref no email%1,            Luis E%Morales-Buenrostro%xref no email%1,            Thierry%Hern\u00e1ndez-Gilsoul%xref no email%1]"
$ws.Cells.Item(36, 9).Value = "_CROSSREF"
